$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Educational Objectives: insert a new paragraph right after the heading
#    (i.e. right before the "[Outline the educational aspects...]" filler
#    paragraph) containing the new explanatory text.
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("[Outline the educational aspects of your game", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$placeholder = $rng.Paragraphs(1)
$placeholder.Range.InsertParagraphBefore()

$rng = $d.Content
$rng.Find.Execute("[Outline the educational aspects of your game", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$placeholder = $rng.Paragraphs(1)
$newPara = $placeholder.Previous()
$newPara.Range.Text = 'The game inspires the player to make decisions on energy efficient upgrades in the game play to receive added benefits while instilling the importance and benefits of those decisions. Throughout the gameplay we motivate the player to go for energy efficient upgrades by providing them with the advantages they can get. For example, even though the initial cost is high they can make more profit as the gameplay goes on with these energy efficient upgrades. '

# ---------------------------------------------------------------------------
# 2) Phase 0 - Player Authentication/Registration: insert a new paragraph
#    right before the "[Explain how you intend...]" filler paragraph.
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("[Explain how you intend to handle player authentication", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$placeholder = $rng.Paragraphs(1)
$placeholder.Range.InsertParagraphBefore()

$rng = $d.Content
$rng.Find.Execute("[Explain how you intend to handle player authentication", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$placeholder = $rng.Paragraphs(1)
$newPara = $placeholder.Previous()
$newPara.Range.Text = 'We implement a component for handling the API calls for the whole game. This component will be the controller for getting and sending API responses and requests. There will be a component for loading the main menu which will talk with this API controller to get necessary values to authenticate the player, do the authentication and load the main menu.'

# ---------------------------------------------------------------------------
# 3) Phase 1 - Player Profile: insert a new paragraph right before the
#    "[Describe how you will design...]" filler paragraph, and drop the
#    now-obsolete lastRenderedPageBreak that used to sit on that paragraph.
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("[Describe how you will design and implement the player profile", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$placeholder = $rng.Paragraphs(1)
$placeholder.Range.InsertParagraphBefore()

$rng = $d.Content
$rng.Find.Execute("[Describe how you will design and implement the player profile", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$placeholder = $rng.Paragraphs(1)
$newPara = $placeholder.Previous()
$newPara.Range.Text = 'We will create another component for the player profile, and we retrieve already available information through our API controller and save it in our player profile component. In the main menu component, we implement an option for the player to create profile. When he clicks that we access the player profile component and retrieve already available values in a form like interface where there will be empty spaces for the users to fill in the missing details. After filling in the details and pressing submit button the missing values will be updated in the profile component locally, then the API controller will be called to save the profile details through API. If any information is not filled in there will be a prompt to fill in the missing information and after completing the profile the game environment will be visible. '

# Rewriting the placeholder text (identical content) drops the stale
# lastRenderedPageBreak marker that used to precede it.
$rng = $d.Content
$rng.Find.Execute("[Describe how you will design and implement the player profile", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$placeholder = $rng.Paragraphs(1)
$placeholder.Range.Text = "[Describe how you will design and implement the player profile. Specify the personal information to be displayed, how missing information will be handled, and the integration with the mock API.]"

# ---------------------------------------------------------------------------
# 4) Phase 2 heading: "Phase 2 - Questionnaire" -> "Phase 2 <EnDash> Questionnaire"
#    and bump the paragraph-mark run size to match the visible run (28/28).
# ---------------------------------------------------------------------------
$dash = [char]0x2013
$rng = $d.Content
$rng.Find.Execute("Phase 2 - Questionnaire", $false, $false, $false, $false, $false, $true, 1, $false, "Phase 2 $dash Questionnaire", 2)

$rng = $d.Content
$rng.Find.Execute("Phase 2 $dash Questionnaire", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$headingPara = $rng.Paragraphs(1)
$headingPara.Range.Font.Size = 14
$headingPara.Range.Font.SizeBi = 14

# ---------------------------------------------------------------------------
# 5) Phase 2 - Questionnaire: insert a new paragraph right before the
#    "[Outline your approach...]" filler paragraph.
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("[Outline your approach for implementing the questionnaire phase", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$placeholder = $rng.Paragraphs(1)
$placeholder.Range.InsertParagraphBefore()

$rng = $d.Content
$rng.Find.Execute("[Outline your approach for implementing the questionnaire phase", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$placeholder = $rng.Paragraphs(1)
$newPara = $placeholder.Previous()
$newPara.Range.Text = 'For the first time players after the profile completion, the main menu component will prompt the user to click a button for the questionnaire asking the player to attempt. After clicking the button, the main menu component will redirect the player to the questionnaire web application in a separate browser tab and waits for the response with the boost score from the web app. Once the player enters the web application, the backend will retrieve the questions from the database one-by-one with the answer choices and send them to the front end, there it will format it for the user to select the answer. Once the answer is selected it will be sent to the backend, and the backend evaluates and stores the results in database. After completing every question, the backend retrieves the evaluation results from the database and the correct answers with the questions and sends them to the frontend to show the results to the user. While this happens, the boost score of the player also will be calculated in the backend. After the player reviews the questions and clicking the button the webapp will redirect the player to the main menu component with the boost score. To make this communication happen between these two applications, we will use the API controller to handle local HTTP communication. '

Write-Output "edit complete"
